# Daily attendance processing - 2025-12-29 10:35:59
# Reorders the "Recorded By" (column G) values so that any "System"/"system"
# entries are moved to the front of the comma-separated list, while the
# relative order of all other entries (and of the System entries among
# themselves) is preserved.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$lastRow = $used.Rows.Count

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)
    $value = $cell.Value2

    if ($value -eq $null) { continue }
    if ($value -eq "") { continue }
    if ($value.IndexOf(",") -lt 0) { continue }

    $parts = $value.Split(",")
    $systemParts = @()
    $restParts = @()

    foreach ($p in $parts) {
        $trimmed = $p.Trim()
        if ($trimmed.ToLower() -eq "system") {
            $systemParts += $trimmed
        } else {
            $restParts += $trimmed
        }
    }

    if ($systemParts.Length -eq 0) { continue }

    $ordered = $systemParts + $restParts
    $newValue = [string]::Join(", ", $ordered)

    if ($newValue -ne $value) {
        $cell.Value = $newValue
    }
}
